$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.729.92"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "2.672.99"

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.65"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.00"
$ws.Range("E6").Value = "  -1.31%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.621"
$ws.Range("E8").Value = "  +4.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.130"
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("E10").Value = "  -1.10%  "

$ws.Range("E11").Value = "  -2.92%  "

$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.59"
$ws.Range("E13").Value = "  -3.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000200"
$ws.Range("E14").Value = "  -6.36%  "

$ws.Range("D15").Value = "3.153.60"
$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("D16").Value = "65.557.78"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").Value = "2.674.27"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.85"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  -1.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.69"
$ws.Range("E20").Value = "  +2.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.69"
$ws.Range("E21").Value = "  -2.98%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E24").Value = "  +3.31%  "

$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.65"
$ws.Range("E26").Value = "  -1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("E27").Value = "  -4.24%  "

$ws.Range("E28").Value = "  -5.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.10"
$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "526.30"
$ws.Range("E32").Value = "  -2.58%  "

$ws.Range("E33").Value = "  -2.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.53"
$ws.Range("E34").Value = "  -1.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  -0.75%  "

$ws.Range("E36").Value = "  -2.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.49"
$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.69"
$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("E40").Value = "  -4.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.71"
$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.58"
$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("E44").Value = "  -3.40%  "

$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("E46").Value = "  -1.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.23"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("E48").Value = "  -2.10%  "

$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("E50").Value = "  +2.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.24"
$ws.Range("E51").Value = "  +1.14%  "
